$wb = $excel.ActiveWorkbook
$wsDatabase = $wb.Worksheets.Item("Database")
$wsNk = $wb.Worksheets.Item("Nk")

# --- Sheet "Nk": append new logged entries (rows 19-25) ---
$wsNk.Activate()
$wsNk.Range("A19").Value = "Porca madonna laida"
$wsNk.Range("A20").Value = "Abuso di bambini canguro"
$wsNk.Range("A21").Value = "Bocca"
$wsNk.Range("A22").Value = "Stocazzo"
$wsNk.Range("A23").Value = "Cacca"
$wsNk.Range("A24").Value = "Suuucaaaaa"
$wsNk.Range("A25").Value = "Porcoddiomaialebastardonegrobruttoinfame"
$wsNk.Range("A23").Select() | Out-Null

# --- Sheet "Database": append new definition row (row 8) ---
$wsDatabase.Activate()
$wsDatabase.Range("A8").Value = "Cacca"
$wsDatabase.Range("B8").Value = "Escremento -- ZATINI"
$wsDatabase.Range("B8").Select() | Out-Null
